# Add, Update and Delete for Education in profile
#
# The original "Skill" sheet is renamed to "Education" (keeping its
# sheetId / position), and a brand-new "Skill" sheet is inserted right
# after it, repopulated with the data that used to live in the old
# "Skill" sheet.

$wb = $excel.ActiveWorkbook

# 1. Rename the existing "Skill" sheet to "Education" - it keeps its
#    sheetId and worksheet part, only the content changes below.
$education = $wb.Worksheets.Item("Skill")
$education.Name = "Education"

# 2. Insert a brand new sheet right after "Education" and call it "Skill".
$skill = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $education)
$skill.Name = "Skill"

# 3. Replace the Education sheet's contents with the education records.
$education.Cells.Clear()

$education.Range("A1").Value = "CollegeName"
$education.Range("B1").Value = "Country"
$education.Range("C1").Value = "Title"
$education.Range("D1").Value = "Degree"
$education.Range("E1").Value = "YearOfPassing"

$education.Range("A2").Value = "BHU"
$education.Range("B2").Value = "India"
$education.Range("C2").Value = "M.A"
$education.Range("D2").Value = "Archeology"
$education.Range("E2").Value = 2013

$education.Range("A3").Value = "ISM"
$education.Range("B3").Value = "INDIA"
$education.Range("C3").Value = "B.Tech"
$education.Range("D3").Value = "Electrical"
$education.Range("E3").Value = 2011

$education.Range("E1:E3").Font.Underline = $true

# 4. Populate the new Skill sheet with what used to be there.
$skill.Range("A1").Value = "Skill"
$skill.Range("B1").Value = "Level"

$skill.Range("A2").Value = "C#"
$skill.Range("B2").Value = 1
$skill.Range("B2").Font.Underline = $true

$skill.Range("A3").Value = "Java"
$skill.Range("B3").Value = 2
